$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update Id and start/end time ---
$ws.Cells.Item(2, 1).Value2 = 111809606      # A2 Id
$ws.Cells.Item(2, 26).Value2 = "14:33"       # Z2 Starttid
$ws.Cells.Item(2, 28).Value2 = "14:33"       # AB2 Sluttid

# --- Rows 3 and 4: the two records swap places entirely ---
# Columns used in these rows: A..I, K, P..W, Z, AB, AD, AE, AG, AT, AW..AY
# (Y/Startdatum and AA/Slutdatum are identical between the two rows, so they
#  are left untouched to avoid the text "2023-08-31" being re-interpreted as
#  a date value when written back through COM.)
$cols = @(1,2,3,4,5,6,7,8,9,11,16,17,18,19,20,21,22,23,26,28,30,31,33,46,49,50,51)

$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Cells.Item(3, $col).Value2
    $row4[$col] = $ws.Cells.Item(4, $col).Value2
}

foreach ($col in $cols) {
    $ws.Cells.Item(3, $col).Value2 = $row4[$col]
    $ws.Cells.Item(4, $col).Value2 = $row3[$col]
}
